# "Only one skill at a time!" — update castTime values and selection/window view.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Abilities")

# Update castTime (column D) values
$ws.Range("D2").Value = 5
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 4
$ws.Range("D3").Value = 15
$ws.Range("D5").Value = 6
$ws.Range("D6").Value = 15

# Update the active selection on the sheet
$ws.Range("E2").Select()
